$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.001.45'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '2.924.04'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.75'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.06'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.84'
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.74'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = '3.404.27'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '60.955.63'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.72'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = '2.924.05'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '431.27'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.37'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.684'
$ws.Range("E21").Value = '  +2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.09'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.61'
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.97'
$ws.Range("E26").Value = '  +1.15%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +5.02%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.08'
$ws.Range("E31").Value = '  -2.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.48'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").Value = '0.0₃0852'
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.64'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.03'
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.51'
$ws.Range("E42").Value = '  -2.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '374.88'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").Value = '2.711.90'
$ws.Range("E45").Value = '  +2.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '130.72'
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.04'
$ws.Range("E48").Value = '  -4.38%  '
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("E51").Value = '  +2.54%  '
